# Scheduled market-data refresh: updates cached price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) on a handful of rows
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4293.364
$ws.Range("J40").Value = 5571
$ws.Range("L40").Value = 5571
$ws.Range("N40").Value = -5921

$ws.Range("H47").Value = 1921.875
$ws.Range("I47").Value = 1166.6666
$ws.Range("K47").Value = 1166.6666
$ws.Range("M47").Value = -194.6666

$ws.Range("H62").Value = 48601.2
$ws.Range("I62").Value = 40335
$ws.Range("K62").Value = 40335
$ws.Range("M62").Value = -39711

$ws.Range("H65").Value = 48601.2
$ws.Range("I65").Value = 40335
$ws.Range("K65").Value = 201675
$ws.Range("M65").Value = -198555

$ws.Range("H80").Value = 535.75
$ws.Range("I80").Value = 528.44446
$ws.Range("J80").Value = 557.6667
$ws.Range("K80").Value = 1585.33338
$ws.Range("L80").Value = 1673.0001
$ws.Range("M80").Value = -587.33338
$ws.Range("N80").Value = -3669.0001

$ws.Range("H83").Value = 535.75
$ws.Range("I83").Value = 528.44446
$ws.Range("J83").Value = 557.6667
$ws.Range("K83").Value = 4756.00014
$ws.Range("L83").Value = 5019.0003
$ws.Range("M83").Value = 235.9998599999999
$ws.Range("N83").Value = -15003.0003

$ws.Range("H98").Value = 8212.25
$ws.Range("I98").Value = 5949.6665
$ws.Range("K98").Value = 5949.6665
$ws.Range("M98").Value = -4451.6665

$ws.Range("H107").Value = 2490.7273
$ws.Range("I107").Value = 427.42856
$ws.Range("J107").Value = 6101.5
$ws.Range("K107").Value = 427.42856
$ws.Range("L107").Value = 6101.5
$ws.Range("M107").Value = 1492.57144
$ws.Range("N107").Value = -9941.5

$ws.Range("H113").Value = 5260.32
$ws.Range("J113").Value = 5312.4375
$ws.Range("L113").Value = 5312.4375
$ws.Range("N113").Value = -11820.4375

$ws.Range("H122").Value = 8212.25
$ws.Range("I122").Value = 5949.6665
$ws.Range("K122").Value = 17848.9995
$ws.Range("M122").Value = -15398.9995

$ws.Range("H132").Value = 3665.138
$ws.Range("I132").Value = 1643.6154
$ws.Range("K132").Value = 4930.8462
$ws.Range("M132").Value = -2400.8462

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2973.1333
$ws.Range("I2").Value = 2738.6155
$ws.Range("K2").Value = 2738.6155
$ws.Range("M2").Value = -2625.6155

$ws.Range("H32").Value = 4197003
$ws.Range("I32").Value = 701101.75
$ws.Range("K32").Value = 701101.75
$ws.Range("M32").Value = -700814.75

$ws.Range("H86").Value = 30314
$ws.Range("J86").Value = 30314
$ws.Range("L86").Value = 30314
$ws.Range("N86").Value = -32686

$ws.Range("H89").Value = 30314
$ws.Range("J89").Value = 30314
$ws.Range("L89").Value = 90942
$ws.Range("N89").Value = -102798

$ws.Range("H97").Value = 18729.615
$ws.Range("I97").Value = 21089.545
$ws.Range("J97").Value = 5750
$ws.Range("K97").Value = 21089.545
$ws.Range("L97").Value = 5750
$ws.Range("M97").Value = -20593.545
$ws.Range("N97").Value = -6742

$ws.Range("H116").Value = 2973.1333
$ws.Range("I116").Value = 2738.6155
$ws.Range("K116").Value = 2738.6155
$ws.Range("M116").Value = -444.6154999999999

$ws.Range("H131").Value = 80715
$ws.Range("J131").Value = 80715
$ws.Range("L131").Value = 80715
$ws.Range("N131").Value = -90795

$ws.Range("H132").Value = 3218.5
$ws.Range("I132").Value = 3218.5
$ws.Range("K132").Value = 9655.5
$ws.Range("M132").Value = -7125.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2973.1333
$ws.Range("I3").Value = 2738.6155
$ws.Range("K3").Value = 2738.6155
$ws.Range("M3").Value = -2624.6155

$ws.Range("H82").Value = 12909.5

$ws.Range("H85").Value = 12909.5

$ws.Range("H134").Value = 6294.8057
$ws.Range("I134").Value = 4920.4585
$ws.Range("J134").Value = 9043.5
$ws.Range("K134").Value = 14761.3755
$ws.Range("L134").Value = 27130.5
$ws.Range("M134").Value = -12226.3755
$ws.Range("N134").Value = -32200.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2556.1428
$ws.Range("I31").Value = 1816.5555
$ws.Range("K31").Value = 1816.5555
$ws.Range("M31").Value = -1521.5555

$ws.Range("H34").Value = 2556.1428
$ws.Range("I34").Value = 1816.5555
$ws.Range("K34").Value = 1816.5555
$ws.Range("M34").Value = -1614.5555

$ws.Range("H86").Value = 66671492
$ws.Range("I86").Value = 100003760
$ws.Range("J86").Value = 6953
$ws.Range("K86").Value = 100003760
$ws.Range("L86").Value = 6953
$ws.Range("M86").Value = -100002637
$ws.Range("N86").Value = -9199

$ws.Range("H89").Value = 66671492
$ws.Range("I89").Value = 100003760
$ws.Range("J89").Value = 6953
$ws.Range("K89").Value = 500018800
$ws.Range("L89").Value = 34765
$ws.Range("M89").Value = -500013184
$ws.Range("N89").Value = -45997

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 8079.5264
$ws.Range("I56").Value = 8079.5264
$ws.Range("K56").Value = 8079.5264
$ws.Range("M56").Value = -7549.5264

$ws.Range("H102").Value = 8553.23
$ws.Range("J102").Value = 8849.333000000001
$ws.Range("L102").Value = 26547.999
$ws.Range("N102").Value = -31415.999

$ws.Range("H117").Value = 3120.4
$ws.Range("I117").Value = 1069.6
$ws.Range("J117").Value = 3530.56
$ws.Range("K117").Value = 3208.8
$ws.Range("L117").Value = 10591.68
$ws.Range("M117").Value = 233.2000000000003
$ws.Range("N117").Value = -17475.68

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 2138.8064
$ws.Range("J107").Value = 2305.1
$ws.Range("L107").Value = 2305.1
$ws.Range("N107").Value = -6145.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1043.1904
$ws.Range("I16").Value = 1097.25
$ws.Range("J16").Value = 870.2
$ws.Range("K16").Value = 1097.25
$ws.Range("L16").Value = 870.2
$ws.Range("M16").Value = -927.25
$ws.Range("N16").Value = -1210.2

$ws.Range("H46").Value = 1224.8334
$ws.Range("I46").Value = 1224.8334
$ws.Range("K46").Value = 1224.8334
$ws.Range("M46").Value = -1036.8334

$ws.Range("H61").Value = 27784978
$ws.Range("I61").Value = 33340424
$ws.Range("J61").Value = 7748.5
$ws.Range("K61").Value = 33340424
$ws.Range("L61").Value = 7748.5
$ws.Range("M61").Value = -33340222
$ws.Range("N61").Value = -8152.5

$ws.Range("H68").Value = 2323.6667
$ws.Range("I68").Value = 2323.6667
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2323.6667
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -1574.6667

$ws.Range("H71").Value = 2323.6667
$ws.Range("I71").Value = 2323.6667
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 11618.3335
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -7874.333500000001

$ws.Range("H82").Value = 859.2222
$ws.Range("I82").Value = 732.25
$ws.Range("J82").Value = 1875
$ws.Range("K82").Value = 732.25
$ws.Range("L82").Value = 1875
$ws.Range("M82").Value = -371.25
$ws.Range("N82").Value = -2597

$ws.Range("H85").Value = 859.2222
$ws.Range("I85").Value = 732.25
$ws.Range("J85").Value = 1875
$ws.Range("K85").Value = 732.25
$ws.Range("L85").Value = 1875
$ws.Range("M85").Value = 515.75
$ws.Range("N85").Value = -4371

$ws.Range("H113").Value = 27784978
$ws.Range("I113").Value = 33340424
$ws.Range("J113").Value = 7748.5
$ws.Range("K113").Value = 33340424
$ws.Range("L113").Value = 7748.5
$ws.Range("M113").Value = -33338254
$ws.Range("N113").Value = -12088.5

$ws.Range("H122").Value = 4556.3335
$ws.Range("I122").Value = 4182.778
$ws.Range("K122").Value = 12548.334
$ws.Range("M122").Value = -10098.334

$ws.Range("H136").Value = 5998
$ws.Range("I136").Value = 4998.3335
$ws.Range("J136").Value = 7497.5
$ws.Range("K136").Value = 14995.0005
$ws.Range("L136").Value = 22492.5
$ws.Range("M136").Value = -12445.0005
$ws.Range("N136").Value = -27592.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 20076.75
$ws.Range("J70").Value = 20076.75
$ws.Range("L70").Value = 20076.75
$ws.Range("N70").Value = -20706.75

$ws.Range("H73").Value = 20076.75
$ws.Range("J73").Value = 20076.75
$ws.Range("L73").Value = 20076.75
$ws.Range("N73").Value = -22260.75

$ws.Range("H122").Value = 2180.2068
$ws.Range("I122").Value = 1733.625
$ws.Range("J122").Value = 2729.8462
$ws.Range("K122").Value = 5200.875
$ws.Range("L122").Value = 8189.5386
$ws.Range("M122").Value = -2750.875
$ws.Range("N122").Value = -13089.5386

$ws.Range("H136").Value = 6189.879
$ws.Range("I136").Value = 4161.8667
$ws.Range("J136").Value = 26470
$ws.Range("K136").Value = 12485.6001
$ws.Range("L136").Value = 79410
$ws.Range("M136").Value = -9935.6001
$ws.Range("N136").Value = -84510
